# This workbook contains weekly "Damasco" (apricot) price records for the
# "Terminal La Palmera de La Serena" market. The edit inserts two new daily
# records (for Castle Brite / Primera and Castle Brite / Segunda, dated
# 2021-12-07) at the top of the existing data block (row 7), pushing all
# the previously-existing records (old rows 7-20) down by two rows to
# rows 9-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 7, shifting rows 7:20 down to 9:22.
$ws.Rows("7:8").Insert()

$newDate = Get-Date -Year 2021 -Month 12 -Day 7 -Hour 0 -Minute 0 -Second 0

# --- New row 7 ---
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Terminal La Palmera de La Serena"
$ws.Range("C7").Value = "Coquimbo"
$ws.Range("D7").Value = $newDate
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100103
$ws.Range("H7").Value = "Frutos de hueso (carozo)"
$ws.Range("I7").Value = 100103003
$ws.Range("J7").Value = "Damasco"
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22500
$ws.Range("Q7").Value = "`$/caja 18 kilos"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1250
$ws.Range("T7").Value = 18

# --- New row 8 ---
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "Terminal La Palmera de La Serena"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = $newDate
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100103
$ws.Range("H8").Value = "Frutos de hueso (carozo)"
$ws.Range("I8").Value = 100103003
$ws.Range("J8").Value = "Damasco"
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 19000
$ws.Range("P8").Value = 18500
$ws.Range("Q8").Value = "`$/caja 18 kilos"
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 1028
$ws.Range("T8").Value = 18

# Ensure the date cells keep the same number format as the rest of column D.
$ws.Range("D7:D8").NumberFormat = $ws.Range("D9").NumberFormat
